$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = "falling"
$ws.Cells.Item(2, 3).Value = -0.4100122451782226
$ws.Cells.Item(2, 4).Value = 0.3157248497009277
$ws.Cells.Item(2, 5).Value = -0.4111190438270569
$ws.Cells.Item(2, 6).Value = -0.0310014113783836
$ws.Cells.Item(2, 7).Value = 0.0068722339347004
$ws.Cells.Item(2, 8).Value = 0.001527163083665

$ws.Cells.Item(3, 1).Value = 100
$ws.Cells.Item(3, 2).Value = "falling"
$ws.Cells.Item(3, 3).Value = -0.3957743644714355
$ws.Cells.Item(3, 4).Value = 0.3475203514099121
$ws.Cells.Item(3, 5).Value = -0.2738307416439056
$ws.Cells.Item(3, 6).Value = -0.0077885319478809
$ws.Cells.Item(3, 7).Value = -0.0080939643085002
$ws.Cells.Item(3, 8).Value = -0.0400116741657257

$ws.Cells.Item(4, 1).Value = 200
$ws.Cells.Item(4, 2).Value = "falling"
$ws.Cells.Item(4, 3).Value = -0.5445261001586914
$ws.Cells.Item(4, 4).Value = 0.2971320152282715
$ws.Cells.Item(4, 5).Value = -0.1380582749843597
$ws.Cells.Item(4, 6).Value = -0.0332921557128429
$ws.Cells.Item(4, 7).Value = -0.0244346093386411
$ws.Cells.Item(4, 8).Value = -0.0705549344420433

$ws.Cells.Item(5, 1).Value = 300
$ws.Cells.Item(5, 2).Value = "falling"
$ws.Cells.Item(5, 3).Value = -0.0015645027160644
$ws.Cells.Item(5, 4).Value = 0.2662014961242676
$ws.Cells.Item(5, 5).Value = -0.2422323226928711
$ws.Cells.Item(5, 6).Value = 0.0039706239476799
$ws.Cells.Item(5, 7).Value = 0.0195476878434419
$ws.Cells.Item(5, 8).Value = -0.038026362657547

$ws.Cells.Item(6, 1).Value = 400
$ws.Cells.Item(6, 2).Value = "falling"
$ws.Cells.Item(6, 3).Value = -0.0723671913146972
$ws.Cells.Item(6, 4).Value = 0.398813247680664
$ws.Cells.Item(6, 5).Value = 0.0261114835739135
$ws.Cells.Item(6, 6).Value = -0.0345138870179653
$ws.Cells.Item(6, 7).Value = 0.1418734490871429
$ws.Cells.Item(6, 8).Value = 0.0546724386513233

$ws.Cells.Item(7, 1).Value = 500
$ws.Cells.Item(7, 2).Value = "falling"
$ws.Cells.Item(7, 3).Value = 0.4648051261901855
$ws.Cells.Item(7, 4).Value = 0.3905239105224609
$ws.Cells.Item(7, 5).Value = 0.401032954454422
$ws.Cells.Item(7, 6).Value = -0.032223142683506
$ws.Cells.Item(7, 7).Value = 0.3484986126422882
$ws.Cells.Item(7, 8).Value = 0.1937969923019409

$ws.Cells.Item(8, 1).Value = 600
$ws.Cells.Item(8, 2).Value = "falling"
$ws.Cells.Item(8, 3).Value = 0.03253173828125
$ws.Cells.Item(8, 4).Value = -0.1416168212890625
$ws.Cells.Item(8, 5).Value = 1.499733686447144
$ws.Cells.Item(8, 6).Value = -0.1188132911920547
$ws.Cells.Item(8, 7).Value = 0.318566232919693
$ws.Cells.Item(8, 8).Value = 0.0813977941870689

$ws.Cells.Item(9, 1).Value = 700
$ws.Cells.Item(9, 2).Value = "falling"
$ws.Cells.Item(9, 3).Value = -0.7911763191223145
$ws.Cells.Item(9, 4).Value = -0.2393178939819336
$ws.Cells.Item(9, 5).Value = 1.927771806716919
$ws.Cells.Item(9, 6).Value = -0.2278527319431305
$ws.Cells.Item(9, 7).Value = 0.3268128931522369
$ws.Cells.Item(9, 8).Value = -0.1646281778812408

$ws.Cells.Item(10, 1).Value = 800
$ws.Cells.Item(10, 2).Value = "falling"
$ws.Cells.Item(10, 3).Value = -0.6626334190368652
$ws.Cells.Item(10, 4).Value = -0.8235597610473633
$ws.Cells.Item(10, 5).Value = 2.30325984954834
$ws.Cells.Item(10, 6).Value = 0.0331394411623477
$ws.Cells.Item(10, 7).Value = -0.1012509167194366
$ws.Cells.Item(10, 8).Value = -0.09605856239795681

$ws.Cells.Item(11, 1).Value = 900
$ws.Cells.Item(11, 2).Value = "falling"
$ws.Cells.Item(11, 3).Value = -1.474053382873535
$ws.Cells.Item(11, 4).Value = -0.7457327842712402
$ws.Cells.Item(11, 5).Value = 2.726755142211914
$ws.Cells.Item(11, 6).Value = 0.7597636580467224
$ws.Cells.Item(11, 7).Value = -0.0694859251379966
$ws.Cells.Item(11, 8).Value = -0.0531452745199203

$ws.Cells.Item(12, 1).Value = 1000
$ws.Cells.Item(12, 2).Value = "falling"
$ws.Cells.Item(12, 3).Value = -2.78084135055542
$ws.Cells.Item(12, 4).Value = -1.402418613433838
$ws.Cells.Item(12, 5).Value = 2.212388038635254
$ws.Cells.Item(12, 6).Value = -0.3869831264019012
$ws.Cells.Item(12, 7).Value = -3.697414636611938
$ws.Cells.Item(12, 8).Value = -1.992184281349182

$ws.Cells.Item(13, 1).Value = 1100
$ws.Cells.Item(13, 2).Value = "falling"
$ws.Cells.Item(13, 3).Value = -3.69126033782959
$ws.Cells.Item(13, 4).Value = -3.209368467330933
$ws.Cells.Item(13, 5).Value = 1.650809645652771
$ws.Cells.Item(13, 6).Value = -1.806633949279785
$ws.Cells.Item(13, 7).Value = -3.700621604919434
$ws.Cells.Item(13, 8).Value = -4.564690589904785

$ws.Cells.Item(14, 1).Value = 1200
$ws.Cells.Item(14, 2).Value = "falling"
$ws.Cells.Item(14, 3).Value = -5.679105281829834
$ws.Cells.Item(14, 4).Value = -4.972231864929199
$ws.Cells.Item(14, 5).Value = 0.8967219591140747
$ws.Cells.Item(14, 6).Value = -0.2874121069908142
$ws.Cells.Item(14, 7).Value = -0.7710646390914917
$ws.Cells.Item(14, 8).Value = 1.056185960769653

$ws.Cells.Item(15, 1).Value = 1300
$ws.Cells.Item(15, 2).Value = "falling"
$ws.Cells.Item(15, 3).Value = -8.219106674194336
$ws.Cells.Item(15, 4).Value = 0.5375771522521973
$ws.Cells.Item(15, 5).Value = -0.9215919971466064
$ws.Cells.Item(15, 6).Value = -0.1162171140313148
$ws.Cells.Item(15, 7).Value = 0.4549418985843658
$ws.Cells.Item(15, 8).Value = -0.4101960062980652

$ws.Cells.Item(16, 1).Value = 1400
$ws.Cells.Item(16, 2).Value = "falling"
$ws.Cells.Item(16, 3).Value = 26.41851043701172
$ws.Cells.Item(16, 4).Value = -1.289778709411621
$ws.Cells.Item(16, 5).Value = -10.24878120422363
$ws.Cells.Item(16, 6).Value = -0.4208861589431762
$ws.Cells.Item(16, 7).Value = -0.4367686510086059
$ws.Cells.Item(16, 8).Value = 0.2874121069908142

$ws.Cells.Item(17, 1).Value = 1500
$ws.Cells.Item(17, 2).Value = "falling"
$ws.Cells.Item(17, 3).Value = -7.613844871520996
$ws.Cells.Item(17, 4).Value = 6.058750152587891
$ws.Cells.Item(17, 5).Value = 5.673562526702881
$ws.Cells.Item(17, 6).Value = -0.6397286057472229
$ws.Cells.Item(17, 7).Value = -0.3787364661693573
$ws.Cells.Item(17, 8).Value = 0.7061602473258972

$ws.Cells.Item(18, 1).Value = 1600
$ws.Cells.Item(18, 2).Value = "falling"
$ws.Cells.Item(18, 3).Value = -1.281689643859863
$ws.Cells.Item(18, 4).Value = -2.930487871170044
$ws.Cells.Item(18, 5).Value = -2.143074989318848
$ws.Cells.Item(18, 6).Value = -0.0201585534960031
$ws.Cells.Item(18, 7).Value = -0.1424843221902847
$ws.Cells.Item(18, 8).Value = 0.6214026808738708

$ws.Cells.Item(19, 1).Value = 1700
$ws.Cells.Item(19, 2).Value = "falling"
$ws.Cells.Item(19, 3).Value = -0.9546890258789062
$ws.Cells.Item(19, 4).Value = 0.6401574611663818
$ws.Cells.Item(19, 5).Value = -0.1838119029998779
$ws.Cells.Item(19, 6).Value = -0.06688974797725671
$ws.Cells.Item(19, 7).Value = 0.08704829961061469
$ws.Cells.Item(19, 8).Value = 0.6957755088806152

$ws.Cells.Item(20, 1).Value = 1800
$ws.Cells.Item(20, 2).Value = "falling"
$ws.Cells.Item(20, 3).Value = 0.0202217102050781
$ws.Cells.Item(20, 4).Value = 1.288838624954224
$ws.Cells.Item(20, 5).Value = 0.7679400444030762
$ws.Cells.Item(20, 6).Value = -0.00167987938039
$ws.Cells.Item(20, 7).Value = 0.9256135821342468
$ws.Cells.Item(20, 8).Value = 0.4283692538738251

$ws.Cells.Item(21, 1).Value = 1900
$ws.Cells.Item(21, 2).Value = "falling"
$ws.Cells.Item(21, 3).Value = 1.244831085205078
$ws.Cells.Item(21, 4).Value = -1.71762228012085
$ws.Cells.Item(21, 5).Value = -0.7488219738006592
$ws.Cells.Item(21, 6).Value = -0.6487388610839844
$ws.Cells.Item(21, 7).Value = 0.8439103364944458
$ws.Cells.Item(21, 8).Value = 0.693637490272522

$ws.Cells.Item(22, 1).Value = 2000
$ws.Cells.Item(22, 2).Value = "falling"
$ws.Cells.Item(22, 3).Value = 2.553339958190918
$ws.Cells.Item(22, 4).Value = 0.4797818660736084
$ws.Cells.Item(22, 5).Value = 0.4767866134643554
$ws.Cells.Item(22, 6).Value = -0.0458148941397666
$ws.Cells.Item(22, 7).Value = 0.3787364661693573
$ws.Cells.Item(22, 8).Value = 0.1919644027948379

$ws.Cells.Item(23, 1).Value = 2100
$ws.Cells.Item(23, 2).Value = "falling"
$ws.Cells.Item(23, 3).Value = 2.489582061767578
$ws.Cells.Item(23, 4).Value = 0.4650382995605469
$ws.Cells.Item(23, 5).Value = 2.65714955329895
$ws.Cells.Item(23, 6).Value = 0.8376489877700806
$ws.Cells.Item(23, 7).Value = 1.98836636543274
$ws.Cells.Item(23, 8).Value = -0.266184538602829

$ws.Cells.Item(24, 1).Value = 2200
$ws.Cells.Item(24, 2).Value = "falling"
$ws.Cells.Item(24, 3).Value = 0.2150793075561523
$ws.Cells.Item(24, 4).Value = 1.105870723724365
$ws.Cells.Item(24, 5).Value = -1.261855840682983
$ws.Cells.Item(24, 6).Value = 1.477683067321777
$ws.Cells.Item(24, 7).Value = 1.546558141708374
$ws.Cells.Item(24, 8).Value = -1.000902652740478

$ws.Cells.Item(25, 1).Value = 2300
$ws.Cells.Item(25, 2).Value = "falling"
$ws.Cells.Item(25, 3).Value = -0.7989382743835449
$ws.Cells.Item(25, 4).Value = -0.0554313659667968
$ws.Cells.Item(25, 5).Value = -2.816707372665405
$ws.Cells.Item(25, 6).Value = -0.024892758578062
$ws.Cells.Item(25, 7).Value = -0.0316122770309448
$ws.Cells.Item(25, 8).Value = -0.06902777403593061

$ws.Cells.Item(26, 1).Value = 2400
$ws.Cells.Item(26, 2).Value = "falling"
$ws.Cells.Item(26, 3).Value = 1.732457160949707
$ws.Cells.Item(26, 4).Value = 0.4553084373474121
$ws.Cells.Item(26, 5).Value = 2.453210830688477
$ws.Cells.Item(26, 6).Value = 0.1531744599342346
$ws.Cells.Item(26, 7).Value = 0.9390525817871094
$ws.Cells.Item(26, 8).Value = -0.3572034537792206

$ws.Cells.Item(27, 1).Value = 2500
$ws.Cells.Item(27, 2).Value = "falling"
$ws.Cells.Item(27, 3).Value = -0.187225341796875
$ws.Cells.Item(27, 4).Value = 0.2255609035491943
$ws.Cells.Item(27, 5).Value = -0.6162976026535034
$ws.Cells.Item(27, 6).Value = 0.1224784851074218
$ws.Cells.Item(27, 7).Value = -0.0070249503478407
$ws.Cells.Item(27, 8).Value = -0.1852448880672454

$ws.Cells.Item(28, 1).Value = 2600
$ws.Cells.Item(28, 2).Value = "falling"
$ws.Cells.Item(28, 3).Value = -0.9302024841308594
$ws.Cells.Item(28, 4).Value = 0.2106423377990722
$ws.Cells.Item(28, 5).Value = -0.1286094188690185
$ws.Cells.Item(28, 6).Value = -0.0273362193256616
$ws.Cells.Item(28, 7).Value = -0.1533271819353103
$ws.Cells.Item(28, 8).Value = -0.1511891484260559

$ws.Cells.Item(29, 1).Value = 2700
$ws.Cells.Item(29, 2).Value = "falling"
$ws.Cells.Item(29, 3).Value = -1.170828819274902
$ws.Cells.Item(29, 4).Value = -0.4288506507873535
$ws.Cells.Item(29, 5).Value = -0.3933718204498291
$ws.Cells.Item(29, 6).Value = -0.0138971842825412
$ws.Cells.Item(29, 7).Value = -0.0740674138069152
$ws.Cells.Item(29, 8).Value = -0.052381694316864

$ws.Cells.Item(30, 1).Value = 2800
$ws.Cells.Item(30, 2).Value = "falling"
$ws.Cells.Item(30, 3).Value = -0.678126335144043
$ws.Cells.Item(30, 4).Value = 0.426605224609375
$ws.Cells.Item(30, 5).Value = -0.0270633697509765
$ws.Cells.Item(30, 6).Value = -0.102472648024559
$ws.Cells.Item(30, 7).Value = -0.0630718395113945
$ws.Cells.Item(30, 8).Value = -0.0003054326225537

$ws.Cells.Item(31, 1).Value = 2900
$ws.Cells.Item(31, 2).Value = "falling"
$ws.Cells.Item(31, 3).Value = -0.493565559387207
$ws.Cells.Item(31, 4).Value = 0.1020381450653076
$ws.Cells.Item(31, 5).Value = 0.0103309154510498
$ws.Cells.Item(31, 6).Value = 0.0174096599221229
$ws.Cells.Item(31, 7).Value = 0.0529925599694252
$ws.Cells.Item(31, 8).Value = -0.0123700210824608

